$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "Parent initial measurements" values in column B (rows 2-14)
$ws.Range("B2").Value  = 47
$ws.Range("B4").Value  = 50
$ws.Range("B5").Value  = 58
$ws.Range("B6").Value  = 56
$ws.Range("B7").Value  = 57
$ws.Range("B8").Value  = 44
$ws.Range("B9").Value  = 62
$ws.Range("B10").Value = 48
$ws.Range("B11").Value = 49
$ws.Range("B12").Value = 57
$ws.Range("B13").Value = 53
$ws.Range("B14").Value = 54

$wb.Save()
